$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to Text format so numeric-looking strings (e.g. "586.07")
# are stored as text, matching the source data which uses inlineStr cells.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '63.708.99'
$ws.Range('E2').Value = '  -0.66%  '

# Row 3
$ws.Range('D3').Value = '3.136.48'
$ws.Range('E3').Value = '  -0.41%  '

# Row 4
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.11%  '

# Row 5
$ws.Range('D5').Value = '586.07'
$ws.Range('E5').Value = '  -0.88%  '

# Row 6
$ws.Range('D6').Value = '146.03'
$ws.Range('E6').Value = '  -1.46%  '

# Row 7
$ws.Range('E7').Value = '  -0.02%  '

# Row 8
$ws.Range('D8').Value = '3.135.04'
$ws.Range('E8').Value = '  -0.30%  '

# Row 9
$ws.Range('D9').Value = '0.529'
$ws.Range('E9').Value = '  -1.22%  '

# Row 10
$ws.Range('E10').Value = '  +3.46%  '

# Row 11
$ws.Range('D11').Value = '5.73'
$ws.Range('E11').Value = '  -1.49%  '

# Row 12
$ws.Range('D12').Value = '0.458'
$ws.Range('E12').Value = '  -3.16%  '

# Row 13
$ws.Range('E13').Value = '  -2.11%  '

# Row 14
$ws.Range('D14').Value = '36.85'
$ws.Range('E14').Value = '  +2.23%  '

# Row 15
$ws.Range('E15').Value = '  -1.79%  '

# Row 16
$ws.Range('D16').Value = '3.652.62'
$ws.Range('E16').Value = '  -0.39%  '

# Row 17
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '63.533.63'
$ws.Range('E17').Value = '  -0.67%  '

# Row 18
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.133.11'
$ws.Range('E18').Value = '  -0.23%  '

# Row 19
$ws.Range('D19').Value = '7.07'

# Row 20
$ws.Range('D20').Value = '463.73'
$ws.Range('E20').Value = '  -2.53%  '

# Row 21
$ws.Range('D21').Value = '14.26'
$ws.Range('E21').Value = '  -0.13%  '

# Row 22
$ws.Range('D22').Value = '0.731'
$ws.Range('E22').Value = '  -0.50%  '

# Row 23
$ws.Range('D23').Value = '7.42'
$ws.Range('E23').Value = '  -2.83%  '

# Row 24
$ws.Range('D24').Value = '12.95'
$ws.Range('E24').Value = '  -3.29%  '

# Row 25
$ws.Range('D25').Value = '81.08'
$ws.Range('E25').Value = '  -1.94%  '

# Row 26
$ws.Range('E26').Value = '  -0.78%  '

# Row 27
$ws.Range('E27').Value = '  -0.07%  '

# Row 28
$ws.Range('D28').Value = '9.23'
$ws.Range('E28').Value = '  +5.44%  '

# Row 29
$ws.Range('E29').Value = '  -0.78%  '

# Row 30
$ws.Range('B30').Value = 'ImmutableX'
$ws.Range('C30').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D30').Value = '2.21'
$ws.Range('E30').Value = '  -1.89%  '

# Row 31
$ws.Range('B31').Value = 'FirstDigitalUSD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D31').Value = '0.999'
$ws.Range('E31').Value = '  -0.16%  '

# Row 32
$ws.Range('D32').Value = '6.99'
$ws.Range('E32').Value = '  +0.60%  '

# Row 33
$ws.Range('D33').Value = '26.92'
$ws.Range('E33').Value = '  -1.54%  '

# Row 34
$ws.Range('E34').Value = '  -0.50%  '

# Row 35
$ws.Range('D35').Value = '0.0₃0844'
$ws.Range('E35').Value = '  -5.32%  '

# Row 36
$ws.Range('E36').Value = '  -1.36%  '

# Row 37
$ws.Range('E37').Value = '  -6.49%  '

# Row 38
$ws.Range('E38').Value = '  -3.62%  '

# Row 39
$ws.Range('D39').Value = '6.01'
$ws.Range('E39').Value = '  -2.49%  '

# Row 40
$ws.Range('D40').Value = '51.03'
$ws.Range('E40').Value = '  +0.12%  '

# Row 41
$ws.Range('D41').Value = '438.69'
$ws.Range('E41').Value = '  -2.08%  '

# Row 42
$ws.Range('E42').Value = '  +0.25%  '

# Row 43
$ws.Range('E43').Value = '  -1.04%  '

# Row 44
$ws.Range('D44').Value = '2.904.86'
$ws.Range('E44').Value = '  -1.78%  '

# Row 45
$ws.Range('D45').Value = '0.277'
$ws.Range('E45').Value = '  -2.33%  '

# Row 46
$ws.Range('E46').Value = '  -3.60%  '

# Row 47
$ws.Range('D47').Value = '37.11'
$ws.Range('E47').Value = '  +3.70%  '

# Row 48
$ws.Range('D48').Value = '126.16'
$ws.Range('E48').Value = '  +1.97%  '

# Row 49
$ws.Range('E49').Value = '  +0.01%  '

# Row 50
$ws.Range('E50').Value = '  -1.71%  '

# Row 51
$ws.Range('D51').Value = '24.23'
$ws.Range('E51').Value = '  -3.08%  '
